# Applies updated transition-probability values to Sheet1
# (added more games -> recomputed simulate-game transition matrix;
#  also reflects drafted optimization-logic state G14 becoming reachable)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2090395480225989
$ws.Range("C2").Value = 0.5423728813559322
$ws.Range("J2").Value = 0.008474576271186441
$ws.Range("P2").Value = 0.1299435028248588
$ws.Range("S2").Value = 0.1101694915254237
$ws.Range("B3").Value = 0.005154639175257732
$ws.Range("C3").Value = 0.02061855670103093
$ws.Range("J3").Value = 0.05670103092783505
$ws.Range("P3").Value = 0.7216494845360825
$ws.Range("S3").Value = 0.1958762886597938
$ws.Range("J4").Value = 0.04081632653061224
$ws.Range("O4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.7142857142857143
$ws.Range("S4").Value = 0.2244897959183673
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.6666666666666666
$ws.Range("B6").Value = 0.03829787234042553
$ws.Range("D6").Value = 0.00425531914893617
$ws.Range("F6").Value = 0.05531914893617021
$ws.Range("J6").Value = 0.2936170212765957
$ws.Range("O6").Value = 0.01702127659574468
$ws.Range("Q6").Value = 0.1361702127659574
$ws.Range("R6").Value = 0.09361702127659574
$ws.Range("S6").Value = 0.3617021276595745
$ws.Range("B7").Value = 0.09693877551020408
$ws.Range("D7").Value = 0.01530612244897959
$ws.Range("F7").Value = 0.03061224489795918
$ws.Range("J7").Value = 0.1530612244897959
$ws.Range("O7").Value = 0.01020408163265306
$ws.Range("Q7").Value = 0.1275510204081633
$ws.Range("R7").Value = 0.1173469387755102
$ws.Range("S7").Value = 0.4489795918367347
$ws.Range("B8").Value = 0.09288537549407115
$ws.Range("D8").Value = 0.01778656126482214
$ws.Range("E8").Value = 0.001976284584980237
$ws.Range("F8").Value = 0.06719367588932806
$ws.Range("J8").Value = 0.1027667984189723
$ws.Range("O8").Value = 0.005928853754940711
$ws.Range("Q8").Value = 0.2055335968379447
$ws.Range("R8").Value = 0.08893280632411067
$ws.Range("S8").Value = 0.41699604743083
$ws.Range("B9").Value = 0.06538461538461539
$ws.Range("D9").Value = 0.003846153846153846
$ws.Range("F9").Value = 0.04230769230769231
$ws.Range("J9").Value = 0.1307692307692308
$ws.Range("O9").Value = 0.003846153846153846
$ws.Range("Q9").Value = 0.1769230769230769
$ws.Range("R9").Value = 0.1269230769230769
$ws.Range("S9").Value = 0.45
$ws.Range("B10").Value = 0.1236263736263736
$ws.Range("D10").Value = 0.02472527472527472
$ws.Range("E10").Value = 0.001373626373626374
$ws.Range("F10").Value = 0.0570054945054945
$ws.Range("J10").Value = 0.1208791208791209
$ws.Range("O10").Value = 0.008241758241758242
$ws.Range("Q10").Value = 0.1868131868131868
$ws.Range("R10").Value = 0.08447802197802198
$ws.Range("S10").Value = 0.3928571428571428
$ws.Range("G11").Value = 0.1335311572700297
$ws.Range("J11").Value = 0.1008902077151335
$ws.Range("K11").Value = 0.1869436201780415
$ws.Range("L11").Value = 0.5548961424332344
$ws.Range("S11").Value = 0.02373887240356083
$ws.Range("G12").Value = 0.6954314720812182
$ws.Range("J12").Value = 0.2487309644670051
$ws.Range("L12").Value = 0.03553299492385787
$ws.Range("S12").Value = 0.02030456852791878
$ws.Range("G13").Value = 0.5555555555555556
$ws.Range("J13").Value = 0.4444444444444444
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.02127659574468085
$ws.Range("H15").Value = 0.1574468085106383
$ws.Range("I15").Value = 0.1063829787234043
$ws.Range("J15").Value = 0.374468085106383
$ws.Range("K15").Value = 0.05106382978723404
$ws.Range("M15").Value = 0.008510638297872341
$ws.Range("O15").Value = 0.05106382978723404
$ws.Range("S15").Value = 0.2297872340425532
$ws.Range("F16").Value = 0.009216589861751152
$ws.Range("H16").Value = 0.1797235023041475
$ws.Range("I16").Value = 0.06912442396313365
$ws.Range("J16").Value = 0.4285714285714285
$ws.Range("K16").Value = 0.1290322580645161
$ws.Range("M16").Value = 0.02304147465437788
$ws.Range("O16").Value = 0.06451612903225806
$ws.Range("S16").Value = 0.09677419354838709
$ws.Range("F17").Value = 0.03368421052631579
$ws.Range("H17").Value = 0.1789473684210526
$ws.Range("I17").Value = 0.09263157894736843
$ws.Range("J17").Value = 0.3936842105263158
$ws.Range("K17").Value = 0.09473684210526316
$ws.Range("M17").Value = 0.01052631578947368
$ws.Range("O17").Value = 0.08842105263157894
$ws.Range("S17").Value = 0.1073684210526316
$ws.Range("F18").Value = 0.01659751037344398
$ws.Range("H18").Value = 0.1203319502074689
$ws.Range("I18").Value = 0.1286307053941909
$ws.Range("J18").Value = 0.4398340248962656
$ws.Range("K18").Value = 0.08298755186721991
$ws.Range("M18").Value = 0.008298755186721992
$ws.Range("N18").Value = 0.004149377593360996
$ws.Range("O18").Value = 0.06224066390041494
$ws.Range("S18").Value = 0.1369294605809129
$ws.Range("F19").Value = 0.02149093351242445
$ws.Range("H19").Value = 0.216252518468771
$ws.Range("I19").Value = 0.09738079247817327
$ws.Range("J19").Value = 0.3492276695768973
$ws.Range("K19").Value = 0.1121558092679651
$ws.Range("M19").Value = 0.01611820013431833
$ws.Range("N19").Value = 0.000671591672263264
$ws.Range("O19").Value = 0.06917394224311618
$ws.Range("S19").Value = 0.1175285426460712
